# CSPro "Tables used to create pngs for export hgts.xlsx" update
#
# - Adds a new "Multiple Files" row to the "File Assoc. Matrix, CSEntry"
#   sheet (as an allowable "Input File" data source), pushing the
#   "Saved Arrays (.sva)" / "Semicolon-delimited" rows down by one.
# - Clears the now-redundant "Input File" (column D) flags that used to
#   be marked individually on several standard file-type rows, since
#   that capability is now represented by the new "Multiple Files" row.
# - Restores the on-screen cell selections that were active when the
#   workbook was last saved.

$wb  = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)   # "File Assoc. Matrix, CSBatch"
$ws4 = $wb.Worksheets.Item(4)   # "File Assoc. Matrix, CSEntry"

$xlPasteAll     = -4104
$xlPasteFormats = -4122
$xlCenter       = -4108

# ---------------------------------------------------------------
# "File Assoc. Matrix, CSEntry" sheet
# ---------------------------------------------------------------

# Push the "Saved Arrays (.sva)" (row 16) and "Semicolon-delimited"
# (row 17) rows down by one, to make room for the new row 16.
$ws4.Rows("16:17").Copy() | Out-Null
$ws4.Rows("17:18").PasteSpecial($xlPasteAll) | Out-Null
$excel.CutCopyMode = 0

# Rebuild row 16 with the correct formatting, copied from cells
# elsewhere on the sheet that already use the desired look.
$ws4.Range("A3").Copy() | Out-Null
$ws4.Range("A16").PasteSpecial($xlPasteFormats) | Out-Null

$ws4.Range("B3").Copy() | Out-Null
$ws4.Range("B16:C16").PasteSpecial($xlPasteFormats) | Out-Null
$ws4.Range("E16:H16").PasteSpecial($xlPasteFormats) | Out-Null

$ws4.Range("C3").Copy() | Out-Null
$ws4.Range("D16").PasteSpecial($xlPasteFormats) | Out-Null

$ws4.Range("H3").Copy() | Out-Null
$ws4.Range("I16").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0

# Two helper cells beyond the printed table (used elsewhere as scratch
# columns) pick up the plain centered look.
$ws4.Range("L16:M16").HorizontalAlignment = $xlCenter

# Fill in the new row's text.
$ws4.Range("A16").Value = "Multiple Files"
$ws4.Range("D16").Value = "Y"

# A handful of rows had their own "Input File" (column D) flag; that's
# now represented by the new "Multiple Files" row instead, so clear it
# (re-using the blank style already used for untouched cells).
$blankTargets = @("D7", "D9", "D10", "D11", "D12", "D14")

$ws4.Range("B3").Copy() | Out-Null
foreach ($cellRef in $blankTargets) {
    $ws4.Range($cellRef).PasteSpecial($xlPasteFormats) | Out-Null
}
$excel.CutCopyMode = 0

foreach ($cellRef in $blankTargets) {
    $ws4.Range($cellRef).ClearContents() | Out-Null
}

# ---------------------------------------------------------------
# Restore view state (active cell / selection) on both sheets
# ---------------------------------------------------------------

$ws3.Activate() | Out-Null
$ws3.Range("A17:XFD17").Select() | Out-Null

$ws4.Activate() | Out-Null
$ws4.Range("O13").Select() | Out-Null
